$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert a new row above the current row 41 ("wohnhaft Sachsen?"), shifting
# the following rows down by one.
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new question/answer pair.
$ws.Range("A41").Value = "juenger als 23 oder vor 1940 geboren?"
$ws.Range("B41").Value = "nein"

# Match formatting used for neighbouring rows: A41 plain (same as A42 below),
# B41 reuses the number-format style from B40 ("Anzahl Kinder"), which is
# what Excel's row insert carries down from the row above.
$ws.Range("B41").NumberFormat = $ws.Range("B40").NumberFormat

# Update the current selection to mirror the saved view state.
$ws.Range("B27").Select()
